$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.189.63"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "2.322.01"
$ws.Range("E3").Value = "  +0.83%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'302.60"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").Value = "'99.55"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.518"
$ws.Range("E9").Value = "  +1.97%  "

$ws.Range("D10").Value = "'36.36"
$ws.Range("E10").Value = "  +5.90%  "

$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("D14").Value = "'6.93"
$ws.Range("E14").Value = "  +1.79%  "

$ws.Range("D15").Value = "2.684.22"
$ws.Range("E15").Value = "  +0.83%  "

$ws.Range("D16").Value = "2.334.19"
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").Value = "'0.800"
$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("D18").Value = "43.095.47"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("D19").Value = "'12.65"
$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'68.28"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").Value = "'239.45"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("E24").Value = "  +2.38%  "

$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'25.53"
$ws.Range("E27").Value = "  +3.52%  "

$ws.Range("D28").Value = "'168.95"
$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("D29").Value = "'34.46"
$ws.Range("E29").Value = "  +1.33%  "

$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("E31").Value = "  -10.55%  "

$ws.Range("E32").Value = "  +3.35%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").Value = "'4.73"
$ws.Range("E34").Value = "  +3.53%  "

$ws.Range("D35").Value = "'17.61"
$ws.Range("E35").Value = "  +2.70%  "

$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").Value = "'0.103"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("E40").Value = "  -2.80%  "

$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("D42").Value = "2.003.76"
$ws.Range("E42").Value = "  +0.27%  "

$ws.Range("E43").Value = "  +1.76%  "

$ws.Range("E44").Value = "  -4.74%  "

$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("D46").Value = "'17.66"
$ws.Range("E46").Value = "  -0.65%  "

$ws.Range("D47").Value = "'2.90"
$ws.Range("E47").Value = "  +0.78%  "

$ws.Range("D48").Value = "'54.95"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("D49").Value = "'75.30"
$ws.Range("E49").Value = "  +7.02%  "

$ws.Range("D50").Value = "2.549.68"
$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("E51").Value = "  +1.73%  "
